# Add two new delivery orders and update order #add6407f... to "Delivered"
# with picked/delivered timestamps, per the uploaded data file change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 (existing order add6407f-8225-46ae-be70-e5a3c9a9b5c7): ---
# status moves from "Pending" to "Delivered", and gets its Picked/Delivered
# timestamps filled in.
$ws.Range("C12").Value = "Delivered"

$ws.Range("D12").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("D12").Value = 45876.96630787037

$ws.Range("E12").NumberFormat = $ws.Range("E2").NumberFormat
$ws.Range("E12").Value = 45876.96645833334

# --- Row 13 (new order): Delivered, with Picked/Delivered timestamps ---
$ws.Range("A13").Value = "89bdc2f6-0e22-47a8-b4f2-b7b5696fc495"
$ws.Range("C13").Value = "Delivered"

$ws.Range("D13").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("D13").Value = 45876.96640046296

$ws.Range("E13").NumberFormat = $ws.Range("E2").NumberFormat
$ws.Range("E13").Value = 45876.96673611111

# --- Row 14 (new order): Pending, no timestamps yet ---
$ws.Range("A14").Value = "0947da20-6ab3-444d-97b4-2aa9c1662a75"
$ws.Range("C14").Value = "Pending"
